$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 was previously the "latest" row (date-only format). Now that a new
# row is appended, it reverts to the normal date-time format used by all
# the other historical rows.
$ws.Range("A41").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new daily row (row 42), whose date cell gets the "latest row"
# date-only format.
$ws.Range("A42").Value = 45627
$ws.Range("A42").NumberFormat = "YYYY-MM-DD"
$ws.Range("B42").Value = 111
$ws.Range("C42").Value = 92
$ws.Range("D42").Value = 100
